$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 235, shifting the existing rows 235:302 down to 236:303
$ws.Rows(235).Insert()

# Populate the newly inserted row 235 with the new weekly price record
$ws.Cells.Item(235, 1).Value  = 11
$ws.Cells.Item(235, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(235, 3).Value  = "Bíobío"
$ws.Cells.Item(235, 4).Value  = 44754
$ws.Cells.Item(235, 5).Value  = 8
$ws.Cells.Item(235, 6).Value  = "Fruta"
$ws.Cells.Item(235, 7).Value  = 100102
$ws.Cells.Item(235, 8).Value  = "Cítricos"
$ws.Cells.Item(235, 9).Value  = 100102005
$ws.Cells.Item(235, 10).Value = "Naranja"
$ws.Cells.Item(235, 11).Value = "Fukumoto"
$ws.Cells.Item(235, 12).Value = "Primera"
$ws.Cells.Item(235, 13).Value = 350
$ws.Cells.Item(235, 14).Value = 7000
$ws.Cells.Item(235, 15).Value = 7500
$ws.Cells.Item(235, 16).Value = 7286
$ws.Cells.Item(235, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(235, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(235, 19).Value = 486
$ws.Cells.Item(235, 20).Value = 15
